# corrected data cleaning for pre/post/total fixation data
#
# 1) Row 1 header cells lose the bold/centered/bordered style (back to the
#    default "Normal" style) and A1's label is cleared.
# 2) Columns C and T (an erroneous "arrg1"/"return2" aggregate column) are
#    cleared out for the data rows (3-8), and the remaining numeric
#    columns D.. are recalculated with corrected totals.
# 3) The stray trailing blank row 10 is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) strip the header row's custom style (bold font + thin box border) ---
$ws.Range("A1:Y1").ClearFormats()
$ws.Range("A1").ClearContents()

# --- 2) clear the retired "arrg1" (C) and "return2" (T) columns ---
$ws.Range("C3:C8").ClearContents()
$ws.Range("T3:T8").ClearContents()

# --- Row 3 : Revisit count ---
$ws.Range("D3").Value = 36
$ws.Range("E3").Value = 20
$ws.Range("F3").Value = 11
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 10
$ws.Range("J3").Value = 6
$ws.Range("Q3").Value = 9
$ws.Range("S3").Value = 1
$ws.Range("U3").Value = 36

# --- Row 4 : Fixation count ---
$ws.Range("D4").Value = 83
$ws.Range("E4").Value = 30
$ws.Range("F4").Value = 18
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 15
$ws.Range("J4").Value = 7
$ws.Range("Q4").Value = 11
$ws.Range("S4").Value = 2
$ws.Range("U4").Value = 127

# --- Row 5 : Dwell time (ms) ---
$ws.Range("D5").Value = 20101.44
$ws.Range("E5").Value = 5905.6
$ws.Range("F5").Value = 5871.9
$ws.Range("G5").Value = 434.38
$ws.Range("H5").Value = 2084.62
$ws.Range("J5").Value = 1635.72
$ws.Range("Q5").Value = 4805.1
$ws.Range("S5").Value = 383.74
$ws.Range("U5").Value = 27774.08

# --- Row 6 : Dwell time (%) ---
$ws.Range("D6").Value = 13.43
$ws.Range("E6").Value = 3.95
$ws.Range("F6").Value = 3.92
$ws.Range("G6").Value = 0.29
$ws.Range("H6").Value = 1.39
$ws.Range("J6").Value = 1.09
$ws.Range("K6").Value = 0.12
$ws.Range("Q6").Value = 3.21
$ws.Range("R6").Value = 0.71
$ws.Range("S6").Value = 0.26
$ws.Range("U6").Value = 18.56
$ws.Range("V6").Value = 0.22
$ws.Range("W6").Value = 0.52
$ws.Range("Y6").Value = 0.12

# --- Row 7 : Fixation duration (ms) ---
$ws.Range("D7").Value = 242.19
$ws.Range("E7").Value = 196.85
$ws.Range("F7").Value = 326.22
$ws.Range("G7").Value = 217.19
$ws.Range("H7").Value = 138.97
$ws.Range("J7").Value = 233.67
$ws.Range("Q7").Value = 436.83
$ws.Range("S7").Value = 191.87
$ws.Range("U7").Value = 218.69

# Row 8 : First fixation duration (ms) - only C8/T8 clearing above applies,
# the remaining values are unchanged.

# --- 3) drop the stray blank row 10 ---
$ws.Rows("10:10").Delete()
